$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.398.17'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.918.45'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.44'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.693'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.66%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '43.98'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.61'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +9.44%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0762'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +3.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0994'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.49'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +9.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.802'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.198.96'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.11'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.915.38'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.424.40'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.81%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0851'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '251.39'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.19'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.69'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.87%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.52'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.79'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.79'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +3.08%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.70%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.37%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.97%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.95'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +6.79%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.33'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +4.65%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0857'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +22.86%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -13.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.858'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.48%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.49%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '106.75'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +10.84%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Gas'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '16.41'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +34.51%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.09'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.341.28'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.28%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0808'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.83%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.44'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.099.66'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.69%  '
